# Weekly update: a new price report row is added for Perejil (Agrícola del
# Norte S.A. de Arica). The new observation becomes the most recent row
# (row 17), pushing the previously-existing rows 17-20 down to 18-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 17, shifting rows 17:20 down to 18:21.
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with this week's data (matches the static columns
# used throughout this subset, plus the new observation's own values).
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C17").Value = "Arica y Parinacota"
$ws.Range("D17").Value = 44438
$ws.Range("E17").Value = 15
$ws.Range("F17").Value = 100112044
$ws.Range("G17").Value = "Perejil"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 300
$ws.Range("K17").Value = 950
$ws.Range("L17").Value = 1000
$ws.Range("M17").Value = 975
$ws.Range("N17").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O17").Value = "Región de Arica y Parinacota"
$ws.Range("P17").Value = 488
$ws.Range("Q17").Value = 2
$ws.Range("R17").Value = "Hortaliza"

# Match the date-cell style used by the other rows in column D.
$ws.Range("D17").NumberFormat = $ws.Range("D18").NumberFormat
